# Apply the "a lot of changes" fixture update to the Man United fixture list.
#
# Logical changes (row numbers as they exist in the original sheet, A=fixture, B=date/time):
#   1. Remove row 1 ("Manchester United v Brighton and Hove Albion " / "18 DecSat12:30")
#      entirely -- the match has been dropped from the list, shifting every
#      following row up by one.
#   2. The "Manchester United v Southampton " fixture's kick-off time changes
#      from "12 FebSat15:00" to "12 FebSat12:30".
#   3. A new fixture is added straight after "Manchester United v Tottenham
#      Hotspur  " / "12 MarSat15:00":
#      "Manchester United v Atletico de Madrid " / "15 MarTue20:00".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Delete the Brighton fixture row (row 1); everything below shifts up.
$ws.Rows.Item(1).Delete()

# 2) Update Southampton's kick-off time. After the deletion above, the
#    Southampton row (originally row 8) is now row 7.
$ws.Cells.Item(7, 2).Value = "12 FebSat12:30"

# 3) Insert a new row right after the Tottenham Hotspur fixture (originally
#    row 11, now row 10 post-deletion) and populate it with the new
#    Atletico de Madrid fixture.
$ws.Rows.Item(11).Insert()
$ws.Cells.Item(11, 1).Value = "Manchester United v Atletico de Madrid "
$ws.Cells.Item(11, 2).Value = "15 MarTue20:00"
